$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(16879.14, 255, 255),
    @(-42.235, 9, 9),
    @(-10.57, 9, 9),
    @(57.229, 9, 9),
    @(0.27, 9, 9),
    @(0.101, 9, 9),
    @(116.591, 9, 9),
    @(-198409.563, 9, 9),
    @(129.059, 9, 9)
)

$startRow = 24
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# The newly appended rows in the source file carry no explicit cell
# style (unlike the existing rows, which use the sheet's column default
# style). Setting the font back to the workbook's base font collapses
# the cell format back to the implicit default style, matching the
# target layout.
$newRange = $ws.Range("A24:C32")
$newRange.Font.Name = "Calibri"
